$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / append opportunity rows (2-14) ---
$ws.Range("A2").Value = '1328113'
$ws.Range("B2").Value = 'https://aiesec.org/opportunity/global-talent/1328113'
$ws.Range("C2").Value = 'Junior Brand Manager - Trainee'
$ws.Range("D2").Value = 'Bruxelles, Belgio'
$ws.Range("E2").Value = 'No'
$ws.Range("F2").Value = '10 applicants'
$ws.Range("G2").Value = '6 - 18 Months'
$ws.Range("H2").Value = 'UCB'

$ws.Range("A3").Value = '1328059'
$ws.Range("B3").Value = 'https://aiesec.org/opportunity/global-talent/1328059'
$ws.Range("C3").Value = 'Marketing Manager'
$ws.Range("D3").Value = 'Delhi, India'
$ws.Range("E3").Value = 'No'
$ws.Range("F3").Value = '2 applicants'
$ws.Range("G3").Value = '6 - 18 Months'
$ws.Range("H3").Value = 'Arvicon International'

$ws.Range("A4").Value = '1327516'
$ws.Range("B4").Value = 'https://aiesec.org/opportunity/global-talent/1327516'
$ws.Range("C4").Value = 'Business development intern'
$ws.Range("D4").Value = 'Sahibzada Ajit Singh Nagar, Punjab, India'
$ws.Range("E4").Value = 'No'
$ws.Range("F4").Value = '1 applicant'
$ws.Range("G4").Value = '9 - 12 Weeks'
$ws.Range("H4").Value = 'AgNext Technologies Private ltd'

$ws.Range("A5").Value = '1327508'
$ws.Range("B5").Value = 'https://aiesec.org/opportunity/global-talent/1327508'
$ws.Range("C5").Value = 'Business Developement Intern'
$ws.Range("D5").Value = 'Sahibzada Ajit Singh Nagar, Punjab, India'
$ws.Range("E5").Value = 'No'
$ws.Range("F5").Value = '1 applicant'
$ws.Range("G5").Value = '9 - 12 Weeks'
$ws.Range("H5").Value = 'AgNext Technologies Private ltd'

$ws.Range("A6").Value = '1326536'
$ws.Range("B6").Value = 'https://aiesec.org/opportunity/global-talent/1326536'
$ws.Range("C6").Value = 'Marketing'
$ws.Range("D6").Value = 'Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E6").Value = 'No'
$ws.Range("F6").Value = '9 applicants'
$ws.Range("G6").Value = '9 - 12 Weeks'
$ws.Range("H6").Value = 'Egypt holiday travel'

$ws.Range("A7").Value = '1326535'
$ws.Range("B7").Value = 'https://aiesec.org/opportunity/global-talent/1326535'
$ws.Range("C7").Value = 'ACCOUNTANT'
$ws.Range("D7").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E7").Value = 'No'
$ws.Range("F7").Value = '21 applicants'
$ws.Range("G7").Value = '9 - 12 Weeks'
$ws.Range("H7").Value = 'Egypt holiday travel'

$ws.Range("A8").Value = '1325846'
$ws.Range("B8").Value = 'https://aiesec.org/opportunity/global-talent/1325846'
$ws.Range("C8").Value = 'Digital Media Strategist - Long Term'
$ws.Range("D8").Value = 'Nugegoda, Sri Lanka'
$ws.Range("E8").Value = 'No'
$ws.Range("F8").Value = '34 applicants'
$ws.Range("G8").Value = '6 - 18 Months'
$ws.Range("H8").Value = 'Brand Corridor (Pvt) Ltd'

$ws.Range("A9").Value = '1325396'
$ws.Range("B9").Value = 'https://aiesec.org/opportunity/global-talent/1325396'
$ws.Range("C9").Value = 'ACE Program | Onboarding & Induction Coordinator'
$ws.Range("D9").Value = 'Budapeste, Hungria'
$ws.Range("E9").Value = 'Yes'
$ws.Range("F9").Value = '214 applicants'
$ws.Range("G9").Value = '6 - 18 Months'
$ws.Range("H9").Value = 'Tata Consultancy Services Ltd.'

$ws.Range("A10").Value = '1325379'
$ws.Range("B10").Value = 'https://aiesec.org/opportunity/global-talent/1325379'
$ws.Range("C10").Value = 'Software Development Intern'
$ws.Range("D10").Value = 'Athens, Greece'
$ws.Range("E10").Value = 'No'
$ws.Range("F10").Value = '105 applicants'
$ws.Range("G10").Value = '9 - 12 Weeks'
$ws.Range("H10").Value = 'Eutopians'

$ws.Range("A11").Value = '1323480'
$ws.Range("B11").Value = 'https://aiesec.org/opportunity/global-talent/1323480'
$ws.Range("C11").Value = 'Mechanical Engineer Intern'
$ws.Range("D11").Value = 'Ankara, Türkiye'
$ws.Range("E11").Value = 'No'
$ws.Range("F11").Value = '50 applicants'
$ws.Range("G11").Value = '9 - 12 Weeks'
$ws.Range("H11").Value = 'AESP MÜHENDİSLİK SANAYİ VE TİCARET ANONİM ŞİRKETİ'

$ws.Range("A12").Value = '1322448'
$ws.Range("B12").Value = 'https://aiesec.org/opportunity/global-talent/1322448'
$ws.Range("C12").Value = 'Web developer'
$ws.Range("D12").Value = 'Giza, El Omraniya, Giza Governorate, Egypt'
$ws.Range("E12").Value = 'No'
$ws.Range("F12").Value = '26 applicants'
$ws.Range("G12").Value = '9 - 12 Weeks'
$ws.Range("H12").Value = 'EG scout shop'

$ws.Range("A13").Value = '1322447'
$ws.Range("B13").Value = 'https://aiesec.org/opportunity/global-talent/1322447'
$ws.Range("C13").Value = 'Social Media&Content Creator'
$ws.Range("D13").Value = 'Giza, El Omraniya, Giza Governorate, Egypt'
$ws.Range("E13").Value = 'No'
$ws.Range("F13").Value = '4 applicants'
$ws.Range("G13").Value = '9 - 12 Weeks'
$ws.Range("H13").Value = 'EG scout shop'

$ws.Range("A14").Value = '1317223'
$ws.Range("B14").Value = 'https://aiesec.org/opportunity/global-talent/1317223'
$ws.Range("C14").Value = 'Accelerate Romania|Software Developer (ONLY EUROPE)'
$ws.Range("D14").Value = 'București, România'
$ws.Range("E14").Value = 'No'
$ws.Range("F14").Value = '112 applicants'
$ws.Range("G14").Value = '9 - 12 Weeks'
$ws.Range("H14").Value = 'Kreston Romania'

# --- Highlight E9 (Premium = Yes) with a yellow fill ---
$ws.Range("E9").Interior.Color = 65535

# --- Resize columns C, D, F, G, H to fit new content ---
$ws.Columns.Item(3).ColumnWidth = 53.17
$ws.Columns.Item(4).ColumnWidth = 69.17
$ws.Columns.Item(6).ColumnWidth = 16.17
$ws.Columns.Item(7).ColumnWidth = 15.17
$ws.Columns.Item(8).ColumnWidth = 51.17
